$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.069367070471715
$ws.Range("D2").Value = 1.052846361169229
$ws.Range("E2").Value = 1.073278608458225
$ws.Range("F2").Value = 1.082737090461309
$ws.Range("I2").Value = 1.051892299978073
$ws.Range("J2").Value = 1.074302157090912
$ws.Range("K2").Value = 1.055593964634548
$ws.Range("L2").Value = 1.075970892234366
$ws.Range("M2").Value = 1.0854045293636
$ws.Range("N2").Value = 1.075827789378161

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.071070588279087
$ws.Range("D3").Value = 1.053687614606996
$ws.Range("E3").Value = 1.074856914062998
$ws.Range("F3").Value = 1.084498441296254
$ws.Range("I3").Value = 1.052385179423964
$ws.Range("J3").Value = 1.075659256580114
$ws.Range("K3").Value = 1.056248072374653
$ws.Range("L3").Value = 1.077364229904162
$ws.Range("M3").Value = 1.086982293062067
$ws.Range("N3").Value = 1.077186816104299

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.072170058252833
$ws.Range("D4").Value = 1.054230578423849
$ws.Range("E4").Value = 1.075875342352951
$ws.Range("F4").Value = 1.085635611722394
$ws.Range("I4").Value = 1.05270164082527
$ws.Range("J4").Value = 1.076534219499119
$ws.Range("K4").Value = 1.056669297561382
$ws.Range("L4").Value = 1.078262479744705
$ws.Range("M4").Value = 1.088000210609618
$ws.Range("N4").Value = 1.078063021571008

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.072631614572007
$ws.Range("D5").Value = 1.054458513326989
$ws.Range("E5").Value = 1.076302822519066
$ws.Range("F5").Value = 1.086113084086075
$ws.Range("I5").Value = 1.052834095423016
$ws.Range("J5").Value = 1.076901306272605
$ws.Range("K5").Value = 1.056845899724633
$ws.Range("L5").Value = 1.078639318273533
$ws.Range("M5").Value = 1.08842743797181
$ws.Range("N5").Value = 1.078430629649819

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.072709073567597
$ws.Range("D6").Value = 1.054496765538061
$ws.Range("E6").Value = 1.076374559578737
$ws.Range("F6").Value = 1.086193219373413
$ws.Range("I6").Value = 1.052856300905166
$ws.Range("J6").Value = 1.076962898249083
$ws.Range("K6").Value = 1.056875523925121
$ws.Range("L6").Value = 1.078702545398452
$ws.Range("M6").Value = 1.08849913035799
$ws.Range("N6").Value = 1.078492309093967

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.072176228168148
$ws.Range("D7").Value = 1.054233625380609
$ws.Range("E7").Value = 1.07588105696714
$ws.Range("F7").Value = 1.08564199404619
$ws.Range("I7").Value = 1.052703412985897
$ws.Range("J7").Value = 1.076539127452388
$ws.Range("K7").Value = 1.056671659211907
$ws.Range("L7").Value = 1.078267518150141
$ws.Range("M7").Value = 1.088005921995167
$ws.Range("N7").Value = 1.078067936494134

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.069943374791409
$ws.Range("D8").Value = 1.053130955638297
$ws.Range("E8").Value = 1.073812600355125
$ws.Range("F8").Value = 1.083332882077527
$ws.Range("I8").Value = 1.052059383765366
$ws.Range("J8").Value = 1.074761458818157
$ws.Range("K8").Value = 1.055815445491883
$ws.Range("L8").Value = 1.076442474660482
$ws.Range("M8").Value = 1.085938371413329
$ws.Range("N8").Value = 1.076287743366541

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.065986562426904
$ws.Range("D9").Value = 1.051177129209921
$ws.Range("E9").Value = 1.070145369123661
$ws.Range("F9").Value = 1.079243823024055
$ws.Range("I9").Value = 1.050905454610062
$ws.Range("J9").Value = 1.071604171129086
$ws.Range("K9").Value = 1.05429097932591
$ws.Range("L9").Value = 1.073200440912244
$ws.Range("M9").Value = 1.082271499842906
$ws.Range("N9").Value = 1.073125971966806

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.063332833442516
$ws.Range("D10").Value = 1.049867065189327
$ws.Range("E10").Value = 1.067684698031412
$ws.Range("F10").Value = 1.076503357811949
$ws.Range("I10").Value = 1.050123078060673
$ws.Range("J10").Value = 1.069481900813786
$ws.Range("K10").Value = 1.053263839170663
$ws.Range("L10").Value = 1.071020783694938
$ws.Range("M10").Value = 1.079810221834258
$ws.Range("N10").Value = 1.071000687784231

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.062179777389083
$ws.Range("D11").Value = 1.049297951054669
$ws.Range("E11").Value = 1.066615255989229
$ws.Range("F11").Value = 1.075313081842417
$ws.Range("I11").Value = 1.049781134562833
$ws.Range("J11").Value = 1.068558640491132
$ws.Range("K11").Value = 1.052816446968763
$ws.Range("L11").Value = 1.07007245468698
$ws.Range("M11").Value = 1.078740312801012
$ws.Range("N11").Value = 1.070076116326018

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.061750866520089
$ws.Range("D12").Value = 1.049086273843678
$ws.Range("E12").Value = 1.066217407349485
$ws.Range("F12").Value = 1.074870396460832
$ws.Range("I12").Value = 1.049653640041012
$ws.Range("J12").Value = 1.068215039747553
$ws.Range("K12").Value = 1.052649864578231
$ws.Range("L12").Value = 1.069719508567052
$ws.Range("M12").Value = 1.078342259730862
$ws.Range("N12").Value = 1.069732027629972

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.061842897448081
$ws.Range("D13").Value = 1.049131692210494
$ws.Range("E13").Value = 1.066302775107383
$ws.Range("F13").Value = 1.074965379779009
$ws.Range("I13").Value = 1.049681009923101
$ws.Range("J13").Value = 1.06828877340408
$ws.Range("K13").Value = 1.052685615299008
$ws.Range("L13").Value = 1.069795248377596
$ws.Range("M13").Value = 1.07842767274591
$ws.Range("N13").Value = 1.069805865996748

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.062144336120505
$ws.Range("D14").Value = 1.04928045955187
$ws.Range("E14").Value = 1.066582382238759
$ws.Range("F14").Value = 1.075276500936099
$ws.Range("I14").Value = 1.049770605679051
$ws.Range("J14").Value = 1.068530251908536
$ws.Range("K14").Value = 1.05280268542783
$ws.Range("L14").Value = 1.070043294305403
$ws.Range("M14").Value = 1.078707422785721
$ws.Range("N14").Value = 1.07004768742838

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.062329980480391
$ws.Range("D15").Value = 1.04937208237713
$ws.Range("E15").Value = 1.066754576074417
$ws.Range("F15").Value = 1.07546811769995
$ws.Range("I15").Value = 1.049825744619674
$ws.Range("J15").Value = 1.068678946800321
$ws.Range("K15").Value = 1.05287476290041
$ws.Range("L15").Value = 1.070196031109258
$ws.Range("M15").Value = 1.078879700559489
$ws.Range("N15").Value = 1.070196593483966

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.063409272762528
$ws.Range("D16").Value = 1.04990479604125
$ws.Range("E16").Value = 1.067755588682239
$ws.Range("F16").Value = 1.076582274398273
$ws.Range("I16").Value = 1.050145704497598
$ws.Range("J16").Value = 1.069543082684512
$ws.Range("K16").Value = 1.053293475175622
$ws.Range("L16").Value = 1.071083624598952
$ws.Range("M16").Value = 1.07988113915084
$ws.Range("N16").Value = 1.07106195654023

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.064085208094909
$ws.Range("D17").Value = 1.05023845484842
$ws.Range("E17").Value = 1.068382427264984
$ws.Range("F17").Value = 1.077280169547159
$ws.Range("I17").Value = 1.05034555462376
$ws.Range("J17").Value = 1.070083970188753
$ws.Range("K17").Value = 1.053555413466054
$ws.Range("L17").Value = 1.071639167083634
$ws.Range("M17").Value = 1.080508190068992
$ws.Range("N17").Value = 1.071603612166768

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.064479087505987
$ws.Range("D18").Value = 1.050432894463665
$ws.Range("E18").Value = 1.068747671074733
$ws.Range("F18").Value = 1.077686890311059
$ws.Range("I18").Value = 1.050461818394505
$ws.Range("J18").Value = 1.070399046485356
$ws.Range("K18").Value = 1.053707943885148
$ws.Range("L18").Value = 1.071962770226518
$ws.Range("M18").Value = 1.080873538048646
$ws.Range("N18").Value = 1.071919135907852

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.064613325899027
$ws.Range("D19").Value = 1.050499163327238
$ws.Range("E19").Value = 1.068872145712582
$ws.Range("F19").Value = 1.077825512799177
$ws.Range("I19").Value = 1.050501409690359
$ws.Range("J19").Value = 1.070506409527151
$ws.Range("K19").Value = 1.053759909934848
$ws.Range("L19").Value = 1.072073037136225
$ws.Range("M19").Value = 1.08099804487615
$ws.Range("N19").Value = 1.072026651417481

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.064012726292989
$ws.Range("D20").Value = 1.050202674857543
$ws.Range("E20").Value = 1.068315212893092
$ws.Range("F20").Value = 1.077205328297126
$ws.Range("I20").Value = 1.050324144226796
$ws.Range("J20").Value = 1.070025981027323
$ws.Range("K20").Value = 1.053527336251833
$ws.Range("L20").Value = 1.07157960777503
$ws.Range("M20").Value = 1.080440954953394
$ws.Range("N20").Value = 1.071545540654077

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.062055587049491
$ws.Range("D21").Value = 1.049236659129428
$ws.Range("E21").Value = 1.066500061882986
$ws.Range("F21").Value = 1.075184899217881
$ws.Range("I21").Value = 1.04974423529644
$ws.Range("J21").Value = 1.068459160862761
$ws.Range("K21").Value = 1.052768222305595
$ws.Range("L21").Value = 1.069970270254712
$ws.Range("M21").Value = 1.078625061204233
$ws.Range("N21").Value = 1.069976495425164

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.060821491757347
$ws.Range("D22").Value = 1.048627647364996
$ws.Range("E22").Value = 1.065355265625448
$ws.Range("F22").Value = 1.073911303937698
$ws.Range("I22").Value = 1.049376835141342
$ws.Range("J22").Value = 1.067470208764945
$ws.Range("K22").Value = 1.052288614451337
$ws.Range("L22").Value = 1.068954390172885
$ws.Range("M22").Value = 1.077479618632395
$ws.Range("N22").Value = 1.068986138901936

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.061476052749939
$ws.Range("D23").Value = 1.048950653195272
$ws.Range("E23").Value = 1.065962484611618
$ws.Range("F23").Value = 1.074586776888002
$ws.Range("I23").Value = 1.049571867019757
$ws.Range("J23").Value = 1.067994838824317
$ws.Range("K23").Value = 1.052543085616562
$ws.Range("L23").Value = 1.069493314055253
$ws.Range("M23").Value = 1.078087197338162
$ws.Range("N23").Value = 1.069511513996168

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.064045478880179
$ws.Range("D24").Value = 1.050218842845199
$ws.Range("E24").Value = 1.068345585350166
$ws.Range("F24").Value = 1.077239146916753
$ws.Range("I24").Value = 1.050333819607319
$ws.Range("J24").Value = 1.070052185113049
$ws.Range("K24").Value = 1.053540023920514
$ws.Range("L24").Value = 1.07160652140802
$ws.Range("M24").Value = 1.080471336844456
$ws.Range("N24").Value = 1.071571781952611

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.0670122209885
$ws.Range("D25").Value = 1.051683544793807
$ws.Range("E25").Value = 1.071096168954119
$ws.Range("F25").Value = 1.080303418477465
$ws.Range("I25").Value = 1.051206060195356
$ws.Range("J25").Value = 1.072423420485558
$ws.Range("K25").Value = 1.054686979248363
$ws.Range("L25").Value = 1.074041753096982
$ws.Range("M25").Value = 1.08322235481133
$ws.Range("N25").Value = 1.073946384751334
